$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.047.60"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.924.06"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.39"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4585"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3818"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07750"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9792"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.64"
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("D12").Value = "1.951.43"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.701"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.960"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07004"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.88"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009499"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.68"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "29.057.95"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.353"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.03"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "2.157.83"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.062"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.98"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.601"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.84"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.831"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09325"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8603"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.096"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05689"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.150"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.004"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02042"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.094"
$ws.Range("E40").Value = "  +12.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.433"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5500"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1755"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.337"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002820"
$ws.Range("E45").Value = "  +7.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.184"
$ws.Range("E46").Value = "  +3.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5172"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.24"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06935"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.57"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.759"
$ws.Range("E51").Value = "  -1.06%  "
